# Add a new "StatQuery" column into the startup sheet and wire up a new
# stats MATCH query next to the existing trial query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing "dbExcel"
# and "WebExcel" filename columns one place to the right (B->C, C->D).
$ws.Columns("B").Insert()

# New column B header + query text.
$ws.Range("B1").Value = "StatQuery"
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN ['Adenocarcinoma of the rectum'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"

# Match the wrap-text formatting used by A2.
$ws.Range("B2").WrapText = $true

# Columns A, C and D keep their original widths automatically (inserting a
# blank column preserves the surviving <col> entries untouched). Column B
# is brand new, so give it the same width as column A.
$ws.Columns("B").ColumnWidth = 75.81640625

$ws.Range("A2").Select()
